# Generate Report for Archive
# The localization status changed from "Ready for handoff" to "In Translation"
# for the two tracked files, across the Overview sheet (columns E/F, one per
# language) and each per-language detail sheet (Status column C). Updating the
# cell text makes the "Status" column content shorter, so Excel auto-fits the
# affected columns back down to the new content width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()

# --- Per-language detail sheets: Status is column C ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus
    $ws.Columns.Item(3).AutoFit()
}
